# Atualização de bases das ligas, do dia: 02-05-2024 às 20:28
#
# The underlying league-data refresh re-sorted some same-date fixtures.
# Column A (the running sequence number) and each row's position on the
# sheet stay exactly where they were; what changes is which fixture's
# data (columns B through AB: id, Div, Date, teams, scores, odds, ...)
# sits in that row. Net effect: the B:AB payloads of these row pairs
# trade places.
#
# Affected row pairs (1-based worksheet rows): (64,66) (65,67) (163,164) (226,227)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowPairs = @(
    @(64, 66),
    @(65, 67),
    @(163, 164),
    @(226, 227)
)

foreach ($pair in $rowPairs) {
    $row1 = $pair[0]
    $row2 = $pair[1]

    $addr1 = "B" + $row1 + ":AB" + $row1
    $addr2 = "B" + $row2 + ":AB" + $row2

    $range1 = $ws.Range($addr1)
    $range2 = $ws.Range($addr2)

    $data1 = $range1.Value2
    $data2 = $range2.Value2

    $range1.Value2 = $data2
    $range2.Value2 = $data1
}
